$d = $word.ActiveDocument

# Find the paragraph whose text is "Linux Bash" and shorten it to "Linux",
# then insert a brand-new paragraph "Linux Mint" right after it (before
# the "Ubuntu" paragraph), matching the style/formatting of its neighbors.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Linux Bash`r") {
        $p.Range.Text = "Linux"
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "Linux Mint"
        break
    }
}
